# Update cryptocurrency price (column D) and 1h volume change (column E)
# values to reflect the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.724.25"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "'1.601.28"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "'211.84"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").Value = "'19.72"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "'0.0845"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "'1.825.92"
$ws.Range("D13").Value = "'1.584.38"
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "'65.00"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "'26.690.25"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "'210.35"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "'7.19"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  -2.64%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'144.06"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").Value = "'7.08"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("D33").Value = "'2.98"
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("D34").Value = "'1.293.18"
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("E37").Value = "  -4.06%  "
$ws.Range("D38").Value = "'1.15"
$ws.Range("E38").Value = "  +7.85%  "
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("E41").Value = "  -2.08%  "
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "'63.02"
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("D45").Value = "'1.738.55"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "'90.53"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("E47").Value = "  -2.42%  "
$ws.Range("D48").Value = "'0.101"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "'0.0516"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").Value = "'7.43"
$ws.Range("E51").Value = "  +0.40%  "
